$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "299猜数字游戏"
$ws.Range("B6").Value = "e"
$ws.Range("C6").Value = "是"
$ws.Range("D6").Value = "和 面试题 16.15. 珠玑妙算 一模一样，不写了"

$ws.Range("D6").Select()
